# default-resume-template.docx edit:
# Insert a new, empty paragraph between the "{/mainSkills}" paragraph and
# the following (empty) section-break paragraph that starts the
# two-column "continuous" section. The new paragraph reuses the
# surrounding paragraph's formatting (pStyle "Normal", the tab stops,
# bidi=0, left alignment) and contains a single empty run - exactly like
# a manual Enter keypress at the end of the "{/mainSkills}" line would
# produce in Word.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("{/mainSkills}", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to a zero-length range right after the matched text, i.e.
    # immediately before the paragraph mark of the "{/mainSkills}"
    # paragraph, and split the paragraph there. The new paragraph that
    # appears before the old paragraph mark inherits that paragraph's
    # formatting (pStyle/tabs/bidi/jc) and is left with a single empty
    # run, matching the target markup.
    $insertionPoint = $d.Range($rng.End, $rng.End)
    $insertionPoint.InsertParagraphAfter()
}
